$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, reusing the exact style of the neighboring
# header cell (G1) rather than building a brand-new style entry.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add value 0 in H2 for the new "Save" column
$ws.Range("H2").Value = 0
